$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A4").Value = "Tobi"
$ws.Range("B4").Value = "Murder"
$ws.Range("C4").Value = "Kreideleichen und so(Amon hats reingeschrieben, Tobi bitte ergänzen)"

$ws.Range("C4").Select()
